$d = $word.ActiveDocument
$sec = $d.Sections(1)

function Rename-LogoShapesIn($story, $oldAlt, $newName) {
    if (-not $story.Exists) { return }
    $paraCount = $story.Range.Paragraphs.Count
    for ($i = 1; $i -le $paraCount; $i++) {
        $p = $story.Range.Paragraphs($i)
        $shapeCount = $p.Range.InlineShapes.Count
        for ($j = 1; $j -le $shapeCount; $j++) {
            $shp = $p.Range.InlineShapes($j)
            if ($shp.AlternativeText -eq $oldAlt) {
                $shp.Name = $newName
            }
        }
    }
}

$btecAlt = "BTec_Logo-Orange"
$pearsonAlt = "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png"

# BTEC logo: appears in the first-page header -> rename image1.jpg to image2.jpg
Rename-LogoShapesIn $sec.Headers(1) $btecAlt "image2.jpg"
Rename-LogoShapesIn $sec.Headers(2) $btecAlt "image2.jpg"
Rename-LogoShapesIn $sec.Headers(3) $btecAlt "image2.jpg"

# Pearson Edexcel logo: appears in both the primary and first-page footers
# -> rename image2.png to image1.png
Rename-LogoShapesIn $sec.Footers(1) $pearsonAlt "image1.png"
Rename-LogoShapesIn $sec.Footers(2) $pearsonAlt "image1.png"
Rename-LogoShapesIn $sec.Footers(3) $pearsonAlt "image1.png"

Write-Host "Renamed header/footer logo shapes"
